# Hub airport selection for COP29
# Reworks the non-direct-flights statistics sheet: the old layout was a
# pandas-style describe() table (row labels "count, mean, std, min, 25%,
# 50%, 75%, max" in column A, 9 metric columns in B:J). The new layout
# drops the describe() row labels and shifts the 9 metric columns one to
# the left so they live in A:I, with fresh per-hub-airport figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old column A describe() labels / percentile style,
#     then drop the now-unused column J entirely -----------------------
$ws.Range("A1:J9").ClearContents()
$ws.Range("A2:A9").ClearFormats()
$ws.Columns.Item(10).Delete()

# --- Header row (A1:I1) ----------------------------------------------------
$headers = @(
    "assumed_load_factor",
    "total_fuel_burn",
    "co2_kg",
    "nox_kg",
    "hc_kg",
    "oc_kg",
    "h2o_kg",
    "so2_kg",
    "sulphates_kg"
)
for ($col = 1; $col -le 9; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Give the now-shifted A1 header cell the same formatting (bold font,
# thin border, centered/top aligned) that the rest of row 1 already has.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# --- Data rows (A2:I9) ------------------------------------------------------
$data = @(
    @(178, 178, 178, 178, 178, 178, 178, 178, 178),
    @(0.6643780106867584, 21371.08718062585, 67511.26440359706, 297.9250240809672, 3.830315180447156, 0.4274217436125169, 26286.43723216979, 25.64530461675102, 0.5233735636071636),
    @(0.07104235951591743, 29143.57529388768, 92064.5543533912, 420.2761367823823, 11.78394463056079, 0.5828715058777537, 35846.59761148185, 34.97229035266522, 0.713720211278882),
    @(0.5245896810666667, 130.2398200270152, 411.427591465341, 1.97183087520901, 0.0129935579494328, 0.002604796400540304, 160.1949786332287, 0.1562877840324182, 0.003189546612906494),
    @(0.6219998006290323, 2670.826547004206, 8437.141061986287, 34.30270701481645, 0.3651244385771767, 0.05341653094008411, 3285.116652815173, 3.204991856405047, 0.06540799706949076),
    @(0.6987738282666667, 7649.937557589101, 24166.15274442397, 86.03013237340585, 0.811637069816892, 0.152998751151782, 9409.423195834595, 9.17992506910692, 0.1873454095736106),
    @(0.7073570327833333, 31964.284773864, 100975.1756006364, 420.0714217200809, 3.582394935444827, 0.6392856954772799, 39316.07027185272, 38.3571417286368, 0.782798810788506),
    @(0.8116680160999999, 112020.975968936, 353874.2630858688, 1589.736681983776, 76.16236322150311, 2.24041951937872, 137785.8004417912, 134.4251711627232, 2.743370840055575)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}

# Fresh export lands the cursor back on A1 (no lingering whole-column
# selection from the old layout).
[void]$ws.Range("A1").Select()

